$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1786743515850144
$ws.Range("C2").Value = 0.579250720461095
$ws.Range("J2").Value = 0.008645533141210375
$ws.Range("P2").Value = 0.1469740634005764
$ws.Range("S2").Value = 0.08645533141210375
$ws.Range("B3").Value = 0.004807692307692308
$ws.Range("C3").Value = 0.02884615384615385
$ws.Range("J3").Value = 0.009615384615384616
$ws.Range("P3").Value = 0.7740384615384616
$ws.Range("S3").Value = 0.1826923076923077
$ws.Range("J4").Value = 0.0576923076923077
$ws.Range("P4").Value = 0.4807692307692308
$ws.Range("S4").Value = 0.4615384615384616
$ws.Range("B6").Value = 0.06024096385542169
$ws.Range("D6").Value = 0.02008032128514056
$ws.Range("E6").Value = 0.004016064257028112
$ws.Range("F6").Value = 0.04417670682730924
$ws.Range("J6").Value = 0.3132530120481928
$ws.Range("O6").Value = 0.04417670682730924
$ws.Range("Q6").Value = 0.1164658634538153
$ws.Range("R6").Value = 0.08032128514056225
$ws.Range("S6").Value = 0.3172690763052209
$ws.Range("B7").Value = 0.1373390557939914
$ws.Range("D7").Value = 0.03433476394849785
$ws.Range("F7").Value = 0.05579399141630902
$ws.Range("J7").Value = 0.09442060085836911
$ws.Range("O7").Value = 0.02575107296137339
$ws.Range("Q7").Value = 0.1502145922746781
$ws.Range("R7").Value = 0.07725321888412018
$ws.Range("S7").Value = 0.4248927038626609
$ws.Range("B8").Value = 0.1318181818181818
$ws.Range("D8").Value = 0.01590909090909091
$ws.Range("F8").Value = 0.05909090909090909
$ws.Range("J8").Value = 0.1181818181818182
$ws.Range("O8").Value = 0.02727272727272727
$ws.Range("Q8").Value = 0.15
$ws.Range("R8").Value = 0.08636363636363636
$ws.Range("S8").Value = 0.4113636363636364
$ws.Range("B9").Value = 0.09146341463414634
$ws.Range("D9").Value = 0.03048780487804878
$ws.Range("F9").Value = 0.07317073170731707
$ws.Range("J9").Value = 0.1036585365853658
$ws.Range("O9").Value = 0.04878048780487805
$ws.Range("Q9").Value = 0.1890243902439024
$ws.Range("R9").Value = 0.06707317073170732
$ws.Range("S9").Value = 0.3963414634146342
$ws.Range("B10").Value = 0.1147540983606557
$ws.Range("D10").Value = 0.02086438152011923
$ws.Range("E10").Value = 0.0007451564828614009
$ws.Range("F10").Value = 0.07526080476900149
$ws.Range("J10").Value = 0.1132637853949329
$ws.Range("O10").Value = 0.02384500745156483
$ws.Range("Q10").Value = 0.2041728763040238
$ws.Range("R10").Value = 0.08047690014903129
$ws.Range("S10").Value = 0.3666169895678092
$ws.Range("F11").Value = 0.002915451895043732
$ws.Range("G11").Value = 0.1516034985422741
$ws.Range("J11").Value = 0.09912536443148688
$ws.Range("K11").Value = 0.2069970845481049
$ws.Range("L11").Value = 0.5306122448979592
$ws.Range("S11").Value = 0.008746355685131196
$ws.Range("G12").Value = 0.7647058823529411
$ws.Range("J12").Value = 0.1764705882352941
$ws.Range("L12").Value = 0.0267379679144385
$ws.Range("S12").Value = 0.03208556149732621
$ws.Range("G13").Value = 0.7796610169491526
$ws.Range("J13").Value = 0.2203389830508475
$ws.Range("F15").Value = 0.02290076335877863
$ws.Range("H15").Value = 0.1641221374045801
$ws.Range("I15").Value = 0.03435114503816794
$ws.Range("J15").Value = 0.3206106870229007
$ws.Range("K15").Value = 0.06106870229007633
$ws.Range("M15").Value = 0.003816793893129771
$ws.Range("N15").Value = 0.003816793893129771
$ws.Range("O15").Value = 0.05725190839694656
$ws.Range("S15").Value = 0.3320610687022901
$ws.Range("F16").Value = 0.013215859030837
$ws.Range("H16").Value = 0.1409691629955947
$ws.Range("I16").Value = 0.07048458149779736
$ws.Range("J16").Value = 0.4581497797356828
$ws.Range("K16").Value = 0.118942731277533
$ws.Range("M16").Value = 0.03524229074889868
$ws.Range("O16").Value = 0.03524229074889868
$ws.Range("S16").Value = 0.1277533039647577
$ws.Range("F17").Value = 0.02546296296296296
$ws.Range("H17").Value = 0.1712962962962963
$ws.Range("I17").Value = 0.09259259259259259
$ws.Range("J17").Value = 0.400462962962963
$ws.Range("K17").Value = 0.1064814814814815
$ws.Range("M17").Value = 0.03703703703703703
$ws.Range("O17").Value = 0.05787037037037037
$ws.Range("S17").Value = 0.1087962962962963
$ws.Range("F18").Value = 0.0267379679144385
$ws.Range("H18").Value = 0.1176470588235294
$ws.Range("I18").Value = 0.0748663101604278
$ws.Range("J18").Value = 0.4705882352941176
$ws.Range("K18").Value = 0.106951871657754
$ws.Range("M18").Value = 0.0213903743315508
$ws.Range("N18").Value = 0.0053475935828877
$ws.Range("O18").Value = 0.06951871657754011
$ws.Range("S18").Value = 0.106951871657754
$ws.Range("F19").Value = 0.02097378277153558
$ws.Range("H19").Value = 0.201498127340824
$ws.Range("I19").Value = 0.06367041198501873
$ws.Range("J19").Value = 0.3737827715355805
$ws.Range("K19").Value = 0.1198501872659176
$ws.Range("M19").Value = 0.02322097378277153
$ws.Range("N19").Value = 0.000749063670411985
$ws.Range("O19").Value = 0.08089887640449438
$ws.Range("S19").Value = 0.1153558052434457
